$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
$sm = $d.SlideMaster
$theme = $sm.Theme
$cs = $theme.ThemeColorScheme
Write-Host "ColorScheme count:" $cs.Count
for ($i=1; $i -le $cs.Count; $i++) {
    $c = $cs.Item($i)
    Write-Host $i ":" $c.RGB
}
